# Anpassung Product Backlog + Sprint PP
#
# The "Aufwandsschätzung" table on Tabelle1 is reorganised:
#   - row 3/4 are trimmed down to a single remaining value each
#   - a brand-new row 5 (A/C/D) and row 6 (A:D, the old row-5 values) appear
#   - row 7/8 become condensed copies (A & C only) of the former row 3/4 data
#   - the "Rene/Dennis/Hans/Vincent" header + SUM formulas move from
#     rows 7/8 down to rows 12/13, with the SUM ranges widened to A2:A8
#   - the grand total in row 16 now sums the new row 13 range
#   - the bar chart series is repointed from Tabelle1!A7:D7 / A8:D8
#     to Tabelle1!A12:D12 / A13:D13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 3: only D3 remains, now holding 2 (was 6) ---------------------
$ws.Cells.Item(3,1).ClearContents()
$ws.Cells.Item(3,2).ClearContents()
$ws.Cells.Item(3,3).ClearContents()
$ws.Cells.Item(3,4).Value = 2

# --- Row 4: only B4 and D4 remain, now holding 4 (were 0.5) ------------
$ws.Cells.Item(4,1).ClearContents()
$ws.Cells.Item(4,2).Value = 4
$ws.Cells.Item(4,3).ClearContents()
$ws.Cells.Item(4,4).Value = 4

# --- Row 5: brand-new row with A5/C5/D5 = 0.5 ---------------------------
$ws.Cells.Item(5,1).Value = 0.5
$ws.Cells.Item(5,2).ClearContents()
$ws.Cells.Item(5,3).Value = 0.5
$ws.Cells.Item(5,4).Value = 0.5

# --- Row 6: brand-new row carrying what used to be row 5 (1,1,1,1) -----
$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = 1
$ws.Cells.Item(6,4).Value = 1

# --- Row 7: condensed copy of the old row-3 data (A & C only) ----------
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).ClearContents()
$ws.Cells.Item(7,3).Value = 6
$ws.Cells.Item(7,4).ClearContents()

# --- Row 8: condensed copy of the old row-4 data (A & C only) ----------
$ws.Cells.Item(8,1).Value = 1.5
$ws.Cells.Item(8,2).ClearContents()
$ws.Cells.Item(8,3).Value = 1.5
$ws.Cells.Item(8,4).ClearContents()

# --- Old header row (7) and SUM row (8) move down to rows 12/13 --------
$ws.Range("A12").Value = $ws.Range("A1").Value2
$ws.Range("B12").Value = $ws.Range("B1").Value2
$ws.Range("C12").Value = $ws.Range("C1").Value2
$ws.Range("D12").Value = $ws.Range("D1").Value2

$ws.Range("A13").Formula = "=SUM(A2:A8)"
$ws.Range("B13:D13").FormulaR1C1 = "=SUM(R[-11]C:R[-5]C)"

# --- Grand total in row 16 now references the relocated SUM row --------
$ws.Range("A16").Formula = "=SUM(A13:D13)"

# --- Re-point the bar-chart series at the relocated header/SUM rows ----
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(,Tabelle1!`$A`$12:`$D`$12,Tabelle1!`$A`$13:`$D`$13,1)"

# --- Cosmetic: leave the active selection on the recomputed grand total
$ws.Range("A16").Select()
